$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BAEPAbCiPC")

# --- "About" sheet: record which state/region this control-lever workbook was copied for ---
# New column C: state name label (B1) + a "last copied/updated" date stamp (C1)
$ws1.Range("B1").Value = "Colorado"
$ws1.Range("C1").Value = 45271
$ws1.Range("C1").NumberFormat = "mm-dd-yy"

# --- "BAEPAbCiPC" sheet: this region/state does not have energy prices affected by ---
# --- production cost changes for any fuel/carrier except electricity (row 2) and heat (row 22) ---
$ws2.Range("B3").Value = 0
$ws2.Range("B4").Value = 0
$ws2.Range("B9").Value = 0
$ws2.Range("B10").Value = 0
$ws2.Range("B11").Value = 0
$ws2.Range("B12").Value = 0
$ws2.Range("B13").Value = 0
$ws2.Range("B14").Value = 0
$ws2.Range("B15").Value = 0
$ws2.Range("B17").Value = 0
$ws2.Range("B18").Value = 0
$ws2.Range("B19").Value = 0
$ws2.Range("B20").Value = 0

# The electricity and heat rows no longer carry the special highlight formatting
$ws2.Range("A2:B2").ClearFormats()
$ws2.Range("A15:B15").ClearFormats()
$ws2.Range("A22:B22").ClearFormats()

# Update the active sheet/selection to match where the editor left off
$ws2.Activate()
$ws2.Range("B18").Select()
